# Add a new BOM line for the VCO control header (10 pin header / 855-M20-9720546)
# right above the existing "20 pin connector" / Jtag adapter block, pushing the
# remaining component rows (and the two summary rows) down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 47, shifting rows 47:53 -> 48:54 and
# the summary rows 56:57 -> 57:58 (formulas/refs auto-adjust).
$ws.Rows(47).Insert()

# Fill in the new component row (order matches how the strings were
# originally added to the workbook: part number, then description, then name).
$ws.Range("C47").Value = "855-M20-9720546"
$ws.Range("B47").Value = "2-row 2.54mm 10 pin header"
$ws.Range("A47").Value = "10 pin header"
$ws.Range("D47").Value = 0.359
$ws.Range("E47").Value = 10
$ws.Range("F47").Formula = "=D47*E47"

# Match the saved selection state from the edit.
[void]$ws.Range("A48").Select()
